$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: make room for 2 new "Bad Drivers" rows (insert before the old Totals row 7) ---
$ws.Range("A7:A8").EntireRow.Insert()

# --- Step 2: make room for 2 new "Good Drivers" rows (insert right after the last existing
#     good-driver row, which after step 1 sits at row 48, i.e. before the trailing blank rows) ---
$ws.Range("A49:A50").EntireRow.Insert()

# --- Bad Drivers data (rows 3-8) ---
$ws.Range("A3").Value = 'Intel(R) Wi-Fi 6 AX200 160MHz - 22.60.0.6'
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 973
$ws.Range("D3").Value = 90.9
$ws.Range("A4").Value = 'Intel(R) Wi-Fi 6 AX200 160MHz - 22.170.0.3'
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 106
$ws.Range("D4").Value = 92.3
$ws.Range("A5").Value = 'Intel(R) Wi-Fi 6E AX210 160MHz - 22.130.0.5'
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 589
$ws.Range("D5").Value = 95.8
$ws.Range("A6").Value = 'Intel(R) Dual Band Wireless-AC 8260 - 20.70.8.1'
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 95
$ws.Range("D6").Value = 97.3
$ws.Range("A7").Value = 'NETGEAR A6100 WiFi Adapter - 1030.25.701.2017'
$ws.Range("B7").Value = 8
$ws.Range("C7").Value = 1485
$ws.Range("D7").Value = 98.6
$ws.Range("A8").Value = 'Intel(R) Wi-Fi 6 AX201 160MHz - 22.250.1.2'
$ws.Range("B8").Value = 17
$ws.Range("C8").Value = 876
$ws.Range("D8").Value = 98.9

# --- Totals row (row 9) ---
$ws.Range("A9").Value = "Totals:"
$ws.Range("B9").Value = 34
$ws.Range("C9").Value = 4124

# --- Good Drivers data (rows 17-50) ---
$ws.Range("A17").Value = 'Intel(R) Wi-Fi 6E AX210 160MHz - 23.120.0.3'
$ws.Range("B17").Value = 34181
$ws.Range("D17").Value = 99.9
$ws.Range("E17").ClearContents()
$ws.Range("A18").Value = 'Intel(R) Wi-Fi 6 AX200 160MHz - 22.160.0.4'
$ws.Range("B18").Value = 96526
$ws.Range("D18").Value = 99.9
$ws.Range("E18").ClearContents()
$ws.Range("A19").Value = 'Intel(R) Wi-Fi 6 AX200 160MHz - 22.230.0.8'
$ws.Range("B19").Value = 328411
$ws.Range("D19").Value = 99.9
$ws.Range("E19").ClearContents()
$ws.Range("A20").Value = 'Intel(R) Wi-Fi 6 AX200 160MHz - 22.200.0.6'
$ws.Range("B20").Value = 143808
$ws.Range("D20").Value = 99.9
$ws.Range("E20").ClearContents()
$ws.Range("A21").Value = 'Intel(R) Wi-Fi 6 AX200 160MHz - 22.190.0.4'
$ws.Range("B21").Value = 287148
$ws.Range("D21").Value = 99.9
$ws.Range("E21").ClearContents()
$ws.Range("A22").Value = 'Intel(R) Wi-Fi 6 AX200 160MHz - 22.250.10.1'
$ws.Range("B22").Value = 69578
$ws.Range("D22").Value = 99.9
$ws.Range("E22").ClearContents()
$ws.Range("A23").Value = 'Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1'
$ws.Range("B23").Value = 56018
$ws.Range("D23").Value = 100
$ws.Range("E23").ClearContents()
$ws.Range("A24").Value = 'Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1'
$ws.Range("B24").Value = 34244
$ws.Range("D24").Value = 100
$ws.Range("E24").ClearContents()
$ws.Range("A25").Value = 'Intel(R) Wi-Fi 6 AX200 160MHz - 22.30.0.11'
$ws.Range("B25").Value = 67111
$ws.Range("D25").Value = 100
$ws.Range("E25").ClearContents()
$ws.Range("A26").Value = 'Intel(R) Wi-Fi 6 AX200 160MHz - 21.30.4.1'
$ws.Range("B26").Value = 13016
$ws.Range("D26").Value = 100
$ws.Range("E26").ClearContents()
$ws.Range("A27").Value = 'Intel(R) Dual Band Wireless-AC 8260 - 22.180.0.4'
$ws.Range("B27").Value = 10456
$ws.Range("D27").Value = 100
$ws.Range("E27").ClearContents()
$ws.Range("A28").Value = 'Intel(R) Dual Band Wireless-AC 8260 - 20.50.0.5'
$ws.Range("B28").Value = 288399
$ws.Range("D28").Value = 100
$ws.Range("E28").ClearContents()
$ws.Range("A29").Value = 'Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4'
$ws.Range("B29").Value = 442178
$ws.Range("D29").Value = 99.9
$ws.Range("E29").Value = '2024-11-10'
$ws.Range("A30").Value = 'Intel(R) Wi-Fi 6E AX210 160MHz - 22.0.1.5'
$ws.Range("B30").Value = 156943
$ws.Range("D30").Value = 100
$ws.Range("E30").Value = '2024-08-13'
$ws.Range("A31").Value = 'Intel(R) Wi-Fi 6 AX200 160MHz - 23.70.2.3'
$ws.Range("B31").Value = 18721
$ws.Range("D31").Value = 99.9
$ws.Range("E31").Value = '2024-07-23'
$ws.Range("A32").Value = 'Intel(R) Wi-Fi 6 AX200 160MHz - 22.10.0.7'
$ws.Range("B32").Value = 66577
$ws.Range("D32").Value = 100
$ws.Range("E32").Value = '2024-05-09'
$ws.Range("A33").Value = 'Intel(R) Wi-Fi 6E AX210 160MHz - 22.130.0.5'
$ws.Range("B33").Value = 18738
$ws.Range("D33").Value = 99.9
$ws.Range("E33").Value = '2024-01-20'
$ws.Range("A34").Value = 'Intel(R) Wi-Fi 6E AX210 160MHz - 23.20.1.1'
$ws.Range("B34").Value = 13533
$ws.Range("D34").Value = 100
$ws.Range("E34").Value = '2023-12-19'
$ws.Range("A35").Value = 'Intel(R) Dual Band Wireless-AC 8260 - 20.70.27.1'
$ws.Range("B35").Value = 17529
$ws.Range("D35").Value = 100
$ws.Range("E35").Value = '2023-09-13'
$ws.Range("A36").Value = 'Intel(R) Wi-Fi 6E AX210 160MHz - 22.170.2.1'
$ws.Range("B36").Value = 19083
$ws.Range("D36").Value = 100
$ws.Range("E36").Value = '2022-11-22'
$ws.Range("A37").Value = 'Intel(R) Wi-Fi 6E AX210 160MHz - 22.100.0.3'
$ws.Range("B37").Value = 12988
$ws.Range("D37").Value = 100
$ws.Range("E37").Value = '2022-05-01'
$ws.Range("A38").Value = 'Intel(R) Wi-Fi 6E AX210 160MHz - 22.110.1.1'
$ws.Range("B38").Value = 42024
$ws.Range("D38").Value = 100
$ws.Range("E38").Value = '2022-05-01'
$ws.Range("A39").Value = 'Intel(R) Dual Band Wireless-AC 8260 - 22.80.1.1'
$ws.Range("B39").Value = 123675
$ws.Range("D39").Value = 100
$ws.Range("E39").Value = '2021-09-11'
$ws.Range("A40").Value = 'Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9'
$ws.Range("B40").Value = 77849
$ws.Range("D40").Value = 99.9
$ws.Range("E40").Value = '2021-08-18'
$ws.Range("A41").Value = 'Intel(R) Wi-Fi 6E AX210 160MHz - 22.70.0.6'
$ws.Range("B41").Value = 15504
$ws.Range("D41").Value = 100
$ws.Range("E41").Value = '2021-06-28'
$ws.Range("A42").Value = 'Intel(R) Wi-Fi 6 AX200 160MHz - 21.60.2.1'
$ws.Range("B42").Value = 26241
$ws.Range("D42").Value = 100
$ws.Range("E42").Value = '2021-01-19'
$ws.Range("A43").Value = 'Intel(R) Wi-Fi 6 AX200 160MHz - 22.0.1.1'
$ws.Range("B43").Value = 15730
$ws.Range("D43").Value = 99.9
$ws.Range("E43").Value = '2020-09-28'
$ws.Range("A44").Value = 'Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2'
$ws.Range("B44").Value = 59673
$ws.Range("D44").Value = 100
$ws.Range("E44").Value = '2020-08-05'
$ws.Range("A45").Value = 'Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6'
$ws.Range("B45").Value = 113652
$ws.Range("D45").Value = 100
$ws.Range("E45").Value = '2019-12-14'
$ws.Range("A46").Value = 'Intel(R) Wi-Fi 6 AX200 160MHz - 21.40.2.2'
$ws.Range("B46").Value = 88435
$ws.Range("D46").Value = 99.9
$ws.Range("E46").Value = '2019-08-31'
$ws.Range("A47").Value = 'Intel(R) Dual Band Wireless-AC 8260 - 20.70.12.5'
$ws.Range("B47").Value = 180575
$ws.Range("D47").Value = 99.9
$ws.Range("E47").Value = '2019-08-25'
$ws.Range("A48").Value = 'Intel(R) Wi-Fi 6 AX200 160MHz - 21.10.1.2'
$ws.Range("B48").Value = 46270
$ws.Range("D48").Value = 100
$ws.Range("E48").Value = '2019-04-23'
$ws.Range("A49").Value = 'Intel(R) Dual Band Wireless-AC 8260 - 20.70.5.2'
$ws.Range("B49").Value = 138724
$ws.Range("D49").Value = 99.9
$ws.Range("E49").Value = '2018-11-25'
$ws.Range("A50").Value = 'Intel(R) Dual Band Wireless-AC 8260 - 20.70.16.4'
$ws.Range("B50").Value = 35023
$ws.Range("D50").Value = 100
$ws.Range("E50").Value = '2018-03-26'
